$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "80÷3=26, 2"
$t.Cell(1, 2).Range.Text  = "15÷8=1, 7"
$t.Cell(1, 3).Range.Text  = "94÷2=47, 0"
$t.Cell(1, 4).Range.Text  = "53÷4=13, 1"
$t.Cell(1, 5).Range.Text  = "36÷8=4, 4"

$t.Cell(5, 1).Range.Text  = "69÷2=34, 1"
$t.Cell(5, 2).Range.Text  = "97÷5=19, 2"
$t.Cell(5, 3).Range.Text  = "34÷2=17, 0"
$t.Cell(5, 4).Range.Text  = "84÷6=14, 0"
$t.Cell(5, 5).Range.Text  = "75÷7=10, 5"

$t.Cell(9, 1).Range.Text  = "15÷9=1, 6"
$t.Cell(9, 2).Range.Text  = "85÷4=21, 1"
$t.Cell(9, 3).Range.Text  = "85÷7=12, 1"
$t.Cell(9, 4).Range.Text  = "63÷7=9, 0"
$t.Cell(9, 5).Range.Text  = "59÷3=19, 2"

$t.Cell(13, 1).Range.Text = "71÷3=23, 2"
$t.Cell(13, 2).Range.Text = "12÷2=6, 0"
$t.Cell(13, 3).Range.Text = "36÷5=7, 1"
$t.Cell(13, 4).Range.Text = "29÷5=5, 4"
$t.Cell(13, 5).Range.Text = "81÷2=40, 1"

$t.Cell(17, 1).Range.Text = "43÷8=5, 3"
$t.Cell(17, 2).Range.Text = "10÷7=1, 3"
$t.Cell(17, 3).Range.Text = "50÷9=5, 5"
$t.Cell(17, 4).Range.Text = "28÷8=3, 4"
$t.Cell(17, 5).Range.Text = "46÷7=6, 4"
